$d = $word.ActiveDocument

# 1. Drop the stray leading ". " from the first bullet's text.
$d.Content.Find.Execute(
    ". One passenger can only be allocated one cabin (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "One passenger can only be allocated one cabin (", 2)

# 2. Add two new bullet points after that paragraph, before the
#    trailing empty paragraph, reusing the same list formatting.
$p1 = $d.Paragraphs.First
$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()

$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Only one guardian is stored for a minor even though they may have more than one onboard the cruise."

$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$p3 = $d.Paragraphs(3)
$p3.Range.Text = "Not every ship needs to be assigned a cruise. For example, a ship that’s currently under maintenance or is being prepared for a cruise but hasn’t been assigned one yet."
